# Add a new row to the end of the log table with the January 6th entry.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRow = $t.Rows.Add()
$idx = $t.Rows.Count

# Match the heading-format flag used by the existing rows so the
# generated <w:trPr> keeps the <w:tblHeader w:val="0"/> element.
$t.Rows.Item($idx).HeadingFormat = 0

$t.Cell($idx, 1).Range.Text = "January 6th"
$t.Cell($idx, 2).Range.Text = "Updated main menu to make buttons more functional, and added text on top of buttons."
$t.Cell($idx, 3).Range.Text = "Updated title page, added name, date, time. Added required files for code to work on github. "
